$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427:492 down to 428:493
$ws.Rows(427).Insert()

# Populate the new row 427 with the new record's data.
# Columns that stay constant throughout the sheet (A,B,C,E,F,G,H,N,Q,R) are
# re-filled the same way they are for every other row in this table.
$ws.Cells.Item(427, 1).Value = 5
$ws.Cells.Item(427, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(427, 3).Value = "Maule"
$ws.Cells.Item(427, 4).Value = 44984
$ws.Cells.Item(427, 5).Value = 7
$ws.Cells.Item(427, 6).Value = 100112023
$ws.Cells.Item(427, 7).Value = "Brócoli"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 5000
$ws.Cells.Item(427, 11).Value = 700
$ws.Cells.Item(427, 12).Value = 700
$ws.Cells.Item(427, 13).Value = 700
$ws.Cells.Item(427, 14).Value = "$/unidad"
$ws.Cells.Item(427, 15).Value = "Región del Maule"
$ws.Cells.Item(427, 16).Value = 700
$ws.Cells.Item(427, 17).Value = 1
$ws.Cells.Item(427, 18).Value = "Hortaliza"
